$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (percent strings, two-dot price strings)
# -- assigned directly as text, no numeric auto-conversion risk.
$textUpdates = @{
    'D2' = '27.373.48'
    'E2' = '  +1.79%  '
    'D3' = '1.826.45'
    'E3' = '  +1.08%  '
    'E4' = '  +0.04%  '
    'E6' = '  +0.00%  '
    'E7' = '  +3.65%  '
    'E8' = '  +3.24%  '
    'E9' = '  +0.96%  '
    'E10' = '  +2.30%  '
    'E11' = '  +0.58%  '
    'D12' = '1.826.72'
    'E12' = '  -2.64%  '
    'E13' = '  +1.56%  '
    'E14' = '  +2.45%  '
    'E15' = '  +0.54%  '
    'E16' = '  +0.38%  '
    'E17' = '  +0.04%  '
    'E18' = '  +0.80%  '
    'E19' = '  +0.01%  '
    'E20' = '  +1.20%  '
    'D21' = '27.373.20'
    'E21' = '  +1.70%  '
    'E22' = '  +3.46%  '
    'E23' = '  +0.99%  '
    'E24' = '  -1.90%  '
    'E25' = '  -0.52%  '
    'E26' = '  +3.13%  '
    'E27' = '  +0.63%  '
    'E28' = '  +2.37%  '
    'E29' = '  +0.37%  '
    'E30' = '  +1.14%  '
    'E31' = '  +6.00%  '
    'E32' = '  +1.49%  '
    'E33' = '  +1.73%  '
    'E34' = '  +0.12%  '
    'E35' = '  +0.03%  '
    'E36' = '  +1.40%  '
    'E37' = '  +0.37%  '
    'E38' = '  +1.06%  '
    'E39' = '  +4.61%  '
    'E40' = '  +0.31%  '
    'E42' = '  +0.96%  '
    'E43' = '  +0.88%  '
    'E44' = '  +2.77%  '
    'E45' = '  -1.45%  '
    'E46' = '  +0.81%  '
    'E47' = '  -0.21%  '
    'E48' = '  +0.92%  '
    'E49' = '  +0.05%  '
    'E50' = '  +0.99%  '
    'E51' = '  +5.28%  '
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Cells whose new values look like plain numbers (e.g. "313.04", "1.0000").
# The source workbook stores these as literal text strings, so force the
# text number format before assigning, then restore the default "Normal"
# style so no stray formatting is left behind on the cell.
$numericTextUpdates = @{
    'D5' = '313.04'
    'D6' = '1.0000'
    'D7' = '0.4602'
    'D8' = '0.3795'
    'D9' = '0.07405'
    'D10' = '0.8773'
    'D13' = '6.711'
    'D14' = '5.436'
    'D15' = '92.99'
    'D16' = '0.07080'
    'D17' = '1.001'
    'D18' = '0.000008806'
    'D19' = '1.0000'
    'D20' = '15.06'
    'D22' = '5.331'
    'D23' = '10.94'
    'D25' = '151.10'
    'D26' = '2.261'
    'D27' = '18.58'
    'D28' = '5.340'
    'D29' = '117.18'
    'D30' = '0.08947'
    'D31' = '0.7972'
    'D32' = '1.194'
    'D33' = '4.546'
    'D34' = '2.938'
    'D35' = '0.9997'
    'D36' = '1.100'
    'D37' = '0.01977'
    'D38' = '0.05256'
    'D39' = '7.354'
    'D40' = '0.5337'
    'D42' = '2.888'
    'D43' = '0.1706'
    'D44' = '8.668'
    'D45' = '0.5092'
    'D46' = '10.62'
    'D47' = '105.18'
    'D48' = '1.684'
    'D49' = '0.9995'
    'D50' = '0.06383'
}
foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
    $cell.Style = "Normal"
}
